# Update "想去人数" (F column) counts on both the "展览" and "全部类型"
# worksheets to reflect the newly generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    3  = 3098
    5  = 2624
    9  = 1368
    11 = 59
    13 = 1179
    14 = 348
    16 = 36
    21 = 2478
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
